$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("품목입출고")
Write-Host $ws.Name
